# Update LR-pairs TPM-derived metrics (ligand/receptor/edge expression &
# specificity values) to reflect the new TPM-based computation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.206015333333333
$ws.Range("H2").Value = 18.618046
$ws.Range("I2").Value = 0.0150172404156507
$ws.Range("J2").Value = 0.0150172404156507
$ws.Range("M2").Value = 0.008446
$ws.Range("N2").Value = 0.025338
$ws.Range("O2").Value = 0.001125187475737063
$ws.Range("P2").Value = 0.001125187475737063
$ws.Range("Q2").Value = 0.05241600550533333
$ws.Range("R2").Value = 0.471744049548
$ws.Range("S2").Value = [double]"1.689721083582261E-05"
$ws.Range("T2").Value = [double]"1.689721083582261E-05"
$ws.Range("G3").Value = 6.206015333333333
$ws.Range("H3").Value = 18.618046
$ws.Range("I3").Value = 0.0150172404156507
$ws.Range("J3").Value = 0.0150172404156507
$ws.Range("N3").Value = 9.970262999999999
$ws.Range("O3").Value = 0.4427506139949732
$ws.Range("P3").Value = 0.4427506139949733
$ws.Range("Q3").Value = 20.625201685122
$ws.Range("R3").Value = 185.626815166098
$ws.Range("S3").Value = 0.006648892414539474
$ws.Range("T3").Value = 0.006648892414539474
$ws.Range("G4").Value = 6.206015333333333
$ws.Range("H4").Value = 18.618046
$ws.Range("I4").Value = 0.0150172404156507
$ws.Range("J4").Value = 0.0150172404156507
$ws.Range("M4").Value = 4.174437666666667
$ws.Range("N4").Value = 12.523313
$ws.Range("O4").Value = 0.5561241985292896
$ws.Range("P4").Value = 0.5561241985292896
$ws.Range("Q4").Value = 25.90662416737756
$ws.Range("R4").Value = 233.159617506398
$ws.Range("S4").Value = 0.008351450790275401
$ws.Range("T4").Value = 0.008351450790275401
$ws.Range("I5").Value = 0.9317452840597572
$ws.Range("J5").Value = 0.9317452840597571
$ws.Range("M5").Value = 0.008446
$ws.Range("N5").Value = 0.025338
$ws.Range("O5").Value = 0.001125187475737063
$ws.Range("P5").Value = 0.001125187475737063
$ws.Range("Q5").Value = 3.252153164435334
$ws.Range("R5").Value = 29.269378479918
$ws.Range("S5").Value = 0.001048388124201111
$ws.Range("T5").Value = 0.00104838812420111
$ws.Range("I6").Value = 0.9317452840597572
$ws.Range("J6").Value = 0.9317452840597571
$ws.Range("N6").Value = 9.970262999999999
$ws.Range("O6").Value = 0.4427506139949732
$ws.Range("P6").Value = 0.4427506139949733
$ws.Range("R6").Value = 11517.22319406909
$ws.Range("S6").Value = 0.4125307966043782
$ws.Range("T6").Value = 0.4125307966043782
$ws.Range("I7").Value = 0.9317452840597572
$ws.Range("J7").Value = 0.9317452840597571
$ws.Range("M7").Value = 4.174437666666667
$ws.Range("N7").Value = 12.523313
$ws.Range("O7").Value = 0.5561241985292896
$ws.Range("P7").Value = 0.5561241985292896
$ws.Range("Q7").Value = 1607.377535802516
$ws.Range("R7").Value = 14466.39782222265
$ws.Range("S7").Value = 0.5181660993311777
$ws.Range("T7").Value = 0.5181660993311777
$ws.Range("G8").Value = 22.00088566666667
$ws.Range("H8").Value = 66.002657
$ws.Range("I8").Value = 0.05323747552459213
$ws.Range("J8").Value = 0.05323747552459213
$ws.Range("M8").Value = 0.008446
$ws.Range("N8").Value = 0.025338
$ws.Range("O8").Value = 0.001125187475737063
$ws.Range("P8").Value = 0.001125187475737063
$ws.Range("Q8").Value = 0.1858194803406667
$ws.Range("R8").Value = 1.672375323066
$ws.Range("S8").Value = [double]"5.990214070012948E-05"
$ws.Range("T8").Value = [double]"5.990214070012948E-05"
$ws.Range("G9").Value = 22.00088566666667
$ws.Range("H9").Value = 66.002657
$ws.Range("I9").Value = 0.05323747552459213
$ws.Range("J9").Value = 0.05323747552459213
$ws.Range("N9").Value = 9.970262999999999
$ws.Range("O9").Value = 0.4427506139949732
$ws.Range("P9").Value = 0.4427506139949733
$ws.Range("Q9").Value = 73.11820544319899
$ws.Range("R9").Value = 658.0638489887909
$ws.Range("S9").Value = 0.02357092497605553
$ws.Range("T9").Value = 0.02357092497605552
$ws.Range("G10").Value = 22.00088566666667
$ws.Range("H10").Value = 66.002657
$ws.Range("I10").Value = 0.05323747552459213
$ws.Range("J10").Value = 0.05323747552459213
$ws.Range("M10").Value = 4.174437666666667
$ws.Range("N10").Value = 12.523313
$ws.Range("O10").Value = 0.5561241985292896
$ws.Range("P10").Value = 0.5561241985292896
$ws.Range("Q10").Value = 91.84132582696012
$ws.Range("R10").Value = 826.571932442641
$ws.Range("S10").Value = 0.02960664840783647
$ws.Range("T10").Value = 0.02960664840783647
